$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item(1)   # "总计"
$wsQ2Src = $wb.Worksheets.Item(2)   # currently named "2022-Q2", holds the old data

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q2" sheet (placed right after the source sheet)
#    and copy the existing quarter's data + formatting into it verbatim,
#    before the source sheet gets renamed/overwritten with Q3 data.
# ---------------------------------------------------------------------
$wsQ2New = $wb.Worksheets.Add($null, $wsQ2Src)
$wsQ2New.Name = "2022-Q2-NEW"

$wsQ2Src.Range("B1:H1").Copy($wsQ2New.Range("B1:H1"))
$wsQ2Src.Range("A2:H3").Copy($wsQ2New.Range("A2:H3"))

# ---------------------------------------------------------------------
# 2) Rename the sheets into their final positions/names:
#      sheet2 (old data, physical file unchanged) -> "2022-Q3"
#      the brand-new sheet (copy of old data)      -> "2022-Q2"
# ---------------------------------------------------------------------
$wsQ2Src.Name = "2022-Q3"
$wsQ2New.Name = "2022-Q2"

# ---------------------------------------------------------------------
# 3) Overwrite the (renamed) "2022-Q3" sheet with the new quarter data.
#    Use an apostrophe prefix so numeric-looking text (fund codes,
#    decimal figures stored as text) is kept as text, not coerced to a
#    number (this also preserves leading zeros, e.g. "014198"), then
#    reset those cells back to the default ("Normal") style, since the
#    apostrophe prefix otherwise tags them with a quote-prefix style.
# ---------------------------------------------------------------------
$wsQ3 = $wsQ2Src

$wsQ3.Range("B1").Value = "基金代码"
$wsQ3.Range("C1").Value = "基金名称"
$wsQ3.Range("D1").Value = "基金规模"
$wsQ3.Range("E1").Value = "股票总仓位"
$wsQ3.Range("F1").Value = "仓位占比"
$wsQ3.Range("G1").Value = "持有市值(亿元)"
$wsQ3.Range("H1").Value = "仓位排名"

$wsQ3.Range("A2").Value = 0
$wsQ3.Range("B2").Value = "'501219"
$wsQ3.Range("C2").Value = "华夏智胜先锋股票（LOF）A"
$wsQ3.Range("D2").Value = "'1.49"
$wsQ3.Range("E2").Value = "'92.63"
$wsQ3.Range("F2").Value = "'0.91"
$wsQ3.Range("G2").Value = "'0.0136"
$wsQ3.Range("H2").Value = 2

$wsQ3.Range("A3").Value = 1
$wsQ3.Range("B3").Value = "'014198"
$wsQ3.Range("C3").Value = "华夏智胜先锋股票（LOF）C"
$wsQ3.Range("D3").Value = "'1.26"
$wsQ3.Range("E3").Value = "'92.63"
$wsQ3.Range("F3").Value = "'0.91"
$wsQ3.Range("G3").Value = "'0.0115"
$wsQ3.Range("H3").Value = 2

# Cells that must end up back on the plain, unstyled ("Normal") style.
$wsQ3.Range("B2").Style = "Normal"
$wsQ3.Range("D2").Style = "Normal"
$wsQ3.Range("E2").Style = "Normal"
$wsQ3.Range("F2").Style = "Normal"
$wsQ3.Range("G2").Style = "Normal"
$wsQ3.Range("B3").Style = "Normal"
$wsQ3.Range("D3").Style = "Normal"
$wsQ3.Range("E3").Style = "Normal"
$wsQ3.Range("F3").Style = "Normal"
$wsQ3.Range("G3").Style = "Normal"

# Re-apply the "总计"-style (bold, centered, bordered - style index 2 in the
# original workbook) to the header row and the index column, matching the
# styling used by the new quarter's data in the target workbook.
$wsTotal.Range("B1").Copy()
$wsQ3.Range("B1:H1").PasteSpecial(-4122)

$wsTotal.Range("A2").Copy()
$wsQ3.Range("A2:A3").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 4) Update the "总计" (total) summary sheet: shift the existing Q2
#    summary row down to row 3, and write the new Q3 summary into row 2.
# ---------------------------------------------------------------------
$wsTotal.Range("A2:D2").Copy($wsTotal.Range("A3:D3"))
$wsTotal.Range("A3").Value = 1

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0.03
